$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Re-order / rewrite column A text labels (rows 2-8) and
# append the new "setup" rows (9-16) used by the image-map tool.
# ---------------------------------------------------------------

$ws.Range("A2").Value = "Bet1"
$ws.Range("A3").Value = "Bet10"
$ws.Range("A4").Value = "Bet100"
$ws.Range("A5").Value = "HitmonSpin"
$ws.Range("A6").Value = "Reel 1"
$ws.Range("A7").Value = "Reel 2"
$ws.Range("A8").Value = "Reel 3"

# ---------------------------------------------------------------
# Row 2 - Bet1
# ---------------------------------------------------------------
$ws.Range("B2").Value = 60
$ws.Range("C2").Value = 60
$ws.Range("D2").Value = 168
$ws.Range("E2").Value = 382

# ---------------------------------------------------------------
# Row 3 - Bet10
# ---------------------------------------------------------------
$ws.Range("B3").Value = 60
$ws.Range("C3").Value = 60
$ws.Range("D3").Formula = "=D2+64+8"
$ws.Range("E3").Value = 382

# ---------------------------------------------------------------
# Row 4 - Bet100
# ---------------------------------------------------------------
$ws.Range("B4").Value = 60
$ws.Range("C4").Value = 60
$ws.Range("D4").Formula = "=D3+64+8"
$ws.Range("E4").Value = 382

# ---------------------------------------------------------------
# Row 5 - HitmonSpin
# ---------------------------------------------------------------
$ws.Range("B5").Value = 70
$ws.Range("C5").Value = 70
$ws.Range("D5").Value = 380
$ws.Range("E5").Value = 375

# ---------------------------------------------------------------
# Row 6 - Reel 1
# ---------------------------------------------------------------
$ws.Range("B6").Value = 67
$ws.Range("C6").Value = 75
$ws.Range("D6").Value = 185
$ws.Range("E6").Value = 210
$ws.Range("F6").Formula = "=FLOOR((D6+((B6-`$B`$9)/2)),1)"
$ws.Range("G6").Formula = "=E6+7"

# ---------------------------------------------------------------
# Row 7 - Reel 2
# ---------------------------------------------------------------
$ws.Range("B7").Value = 67
$ws.Range("C7").Value = 75
$ws.Range("D7").Value = 265
$ws.Range("E7").Value = 210
$ws.Range("F7").Formula = "=FLOOR((D7+((B7-`$B`$9)/2)),1)"
$ws.Range("G7").Formula = "=E7+7"
$ws.Range("I7").Formula = "=F7-F6"

# ---------------------------------------------------------------
# Row 8 - Reel 3
# ---------------------------------------------------------------
$ws.Range("B8").Value = 67
$ws.Range("C8").Value = 75
$ws.Range("D8").Value = 343
$ws.Range("E8").Value = 210
$ws.Range("F8").Formula = "=FLOOR((D8+((B8-`$B`$9)/2)),1)"
$ws.Range("G8").Formula = "=E8+7"
$ws.Range("I8").Formula = "=F8-F7"

# ---------------------------------------------------------------
# New rows 9-16 - symbol setup.
# NOTE: the text labels are written in the same order the
# original author typed them (Raikou, Pikachu, Magikarp,
# Charizard, Articuno, Voltorb, Blank, Jigglypuff) so that the
# shared-string table ends up with the same ordering, even though
# the rows themselves land in a different numeric order.
# ---------------------------------------------------------------
$ws.Range("A15").Value = "Raikou"
$ws.Range("A13").Value = "Pikachu"
$ws.Range("A10").Value = "Magikarp"
$ws.Range("A14").Value = "Charizard"
$ws.Range("A16").Value = "Articuno"
$ws.Range("A11").Value = "Voltorb"
$ws.Range("A9").Value = "Blank"
$ws.Range("A12").Value = "Jigglypuff"

$ws.Range("B9").Value = 50
$ws.Range("C9").Value = 50

$ws.Range("B10").Value = 50
$ws.Range("C10").Value = 50

$ws.Range("B11").Value = 50
$ws.Range("C11").Value = 50

$ws.Range("B12").Value = 50
$ws.Range("C12").Value = 50

$ws.Range("B13").Value = 50
$ws.Range("C13").Value = 50

$ws.Range("B14").Value = 50
$ws.Range("C14").Value = 50

$ws.Range("B15").Value = 50
$ws.Range("C15").Value = 50

$ws.Range("B16").Value = 50
$ws.Range("C16").Value = 50

# ---------------------------------------------------------------
# Styling
# ---------------------------------------------------------------

# Rows 2-4 (Bet buttons): B,C,E use the new "no border center" style
$ws.Range("B2:C2").Style = "BetNumbers"
$ws.Range("E2").Style = "BetNumbers"
$ws.Range("B3:C3").Style = "BetNumbers"
$ws.Range("E3").Style = "BetNumbers"
$ws.Range("B4:C4").Style = "BetNumbers"
$ws.Range("E4").Style = "BetNumbers"

# Row 5 (HitmonSpin): D,E use the new style too
$ws.Range("D5:E5").Style = "BetNumbers"

# New setup rows use default style already (style index 1 = Normal w/ col style)

# Selection
$ws.Range("D5:E5").Select()

# Page setup - orientation portrait
$ws.PageSetup.Orientation = 1
